# Update readings to match table in paper.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C16").Value = "DSC Ch. 17, 18.0 - 18.2, 18.8, 18.10"
$ws.Range("C23").Value = "DSC Ch. 23.0 - 23.4 (skipping 23.2.3 - 23.3)"
$ws.Range("C27").Value = "Junquiera 2011; DSC Ch. 23.6 - 23.9"

# Leave the selection where the author ended up after editing the table.
$ws.Range("C28").Select() | Out-Null
